$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Final Summary Table")
$ws.Range("C1").EntireColumn.Insert()
$ws.Range("C1").ColumnWidth = 7.14
$ws.Range("B1").Value = "table of measured values during/post aerosol exposure (summary per strain)"
$ws.Range("B1:I1").HorizontalAlignment = -4108
$ws.Range("C3:C10").HorizontalAlignment = -4108
$ws.Range("B1:I1").Merge()
$ws.Range("C2").Value = "N"
